# cifar10-on-ResNet18-add-gaussian-blur-gradient-based-attacks.xlsx
# Commit: "how to do it for full cifar10"
#
# Adds a new "spacing" worksheet (as the last tab, made active) that
# tabulates spacing-attack results for the imagenet and cifar10 datasets,
# and nudges the selection on the previously-active "GradientSignRound"
# sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "spacing" sheet after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws8.Name = "spacing"

# ---------------------------------------------------------------------
# 2. First block: "dataset:  imagenet"
# ---------------------------------------------------------------------
$ws8.Range("A1").Value = "dataset:  imagenet"

$ws8.Range("A2").Value = "spacing"
$ws8.Range("B2").Value = " correct"
$ws8.Range("C2").Value = " counter"
$ws8.Range("D2").Value = " correct rate (%)"
$ws8.Range("E2").Value = " time (sec)"

$imagenetData = @(
    @(1,   18, 20, 0.9,  7.4878764152526802),
    @(2,   17, 20, 0.85, 1.8398396968841499),
    @(4,   17, 20, 0.85, 1.8233780860900799),
    @(8,   17, 20, 0.85, 1.83852338790893),
    @(16,  16, 20, 0.8,  1.9674994945526101),
    @(32,  13, 20, 0.65, 1.8854739665985101),
    @(64,  9,  20, 0.45, 1.8965969085693299),
    @(128, 2,  20, 0.1,  1.76455473899841)
)

$row = 3
foreach ($r in $imagenetData) {
    $ws8.Cells.Item($row, 1).Value = $r[0]
    $ws8.Cells.Item($row, 2).Value = $r[1]
    $ws8.Cells.Item($row, 3).Value = $r[2]
    $ws8.Cells.Item($row, 4).Value = $r[3]
    $ws8.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 3. Second block: "dataset:  cifar10"
# ---------------------------------------------------------------------
$ws8.Range("A12").Value = "dataset:  cifar10"

$ws8.Range("A13").Value = "spacing"
$ws8.Range("B13").Value = " correct"
$ws8.Range("C13").Value = " counter"
$ws8.Range("D13").Value = " correct rate (%)"
$ws8.Range("E13").Value = " time (sec)"
$ws8.Range("F13").Value = " sum_difference"

$cifarData = @(
    @(1,   20, 20, 1,    7.52235531806945,    0.00099651515483856201),
    @(2,   20, 20, 1,    2.03090095520019,    119.672306060791),
    @(4,   20, 20, 1,    2.0150196552276598,  240.92734432220399),
    @(8,   20, 20, 1,    2.0048916339874201,  491.96847343444801),
    @(16,  19, 20, 0.95, 2.0036807060241699,  1013.6589012145899),
    @(32,  18, 20, 0.9,  2.0042469501495299,  2207.9646835327098),
    @(64,  16, 20, 0.8,  2.0304780006408598,  5181.68605041503),
    @(128, 7,  20, 0.35, 1.9716362953186,     17603.9998779296)
)

$row = 14
foreach ($r in $cifarData) {
    $ws8.Cells.Item($row, 1).Value = $r[0]
    $ws8.Cells.Item($row, 2).Value = $r[1]
    $ws8.Cells.Item($row, 3).Value = $r[2]
    $ws8.Cells.Item($row, 4).Value = $r[3]
    $ws8.Cells.Item($row, 5).Value = $r[4]
    $ws8.Cells.Item($row, 6).Value = $r[5]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4. Window/selection bookkeeping.
#    - GradientSignRound (previously active/selected tab) keeps "working"
#      but the selection moves further down/right, and it stops being
#      the tab shown when the file is reopened.
#    - The new "spacing" sheet becomes the active tab, with H4 selected.
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("GradientSignRound")
$ws7.Activate()
$ws7.Range("J39").Select()

$ws8.Activate()
$ws8.Range("H4").Select()
